$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,8
$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = -9.370760740737579
$data[0,3] = -3.25683893229624
$data[0,4] = -5.812264925753899
$data[0,5] = 0.7681222558021545
$data[0,6] = -0.0100534334778785
$data[0,7] = -2.038516759872437
$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = -8.167156538812145
$data[1,3] = -0.02954524458807972
$data[1,4] = -3.384951503568034
$data[1,5] = -0.3067295849323272
$data[1,6] = -1.010869383811951
$data[1,7] = -0.750079333782196
$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = -0.9182361909166679
$data[2,3] = -0.4588523148410963
$data[2,4] = -0.8686674579775859
$data[2,5] = -0.4047339260578155
$data[2,6] = -2.28332781791687
$data[2,7] = 0.921320617198944
$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = 5.347650556003387
$data[3,3] = 3.73906264585607
$data[3,4] = -1.996283573262832
$data[3,5] = -0.3051316738128662
$data[3,6] = 1.929260492324829
$data[3,7] = -1.74610161781311
$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = 10.92234480758596
$data[4,3] = 3.812681239115181
$data[4,4] = -1.331586267613695
$data[4,5] = 3.213436841964722
$data[4,6] = 1.757220268249511
$data[4,7] = -4.128565788269043
$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = 9.719275569484333
$data[5,3] = 0.7055869749765016
$data[5,4] = -0.7485139332207429
$data[5,5] = -2.472412109375
$data[5,6] = -2.355233192443848
$data[5,7] = 1.659016251564026
$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = 13.93544812655545
$data[6,3] = -16.89472820316476
$data[6,4] = -4.748856436612479
$data[6,5] = -5.059939384460449
$data[6,6] = -0.0345545150339603
$data[6,7] = 1.632384657859802
$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = 7.450619108536684
$data[7,3] = -0.9767160588260171
$data[7,4] = 19.51188305815956
$data[7,5] = -4.271111488342285
$data[7,6] = 6.537594318389893
$data[7,7] = 0.8408931493759155
$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = 9.203229472647893
$data[8,3] = -3.916383311759566
$data[8,4] = 12.78814014158619
$data[8,5] = 0.7947538495063782
$data[8,6] = -4.166182518005371
$data[8,7] = -1.592170834541321
$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = 5.404519542849628
$data[9,3] = -19.64550725151458
$data[9,4] = -2.510032804843359
$data[9,5] = 4.600943565368652
$data[9,6] = -7.544669151306152
$data[9,7] = -3.569834232330322
$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = 13.85390501755906
$data[10,3] = 10.65629918758675
$data[10,4] = -23.50188240638136
$data[10,5] = 4.194012641906738
$data[10,6] = 4.047538757324219
$data[10,7] = -7.187472343444824
$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = -15.78896728791737
$data[11,3] = -2.495026884036708
$data[11,4] = -9.361624100628198
$data[11,5] = -5.213337898254395
$data[11,6] = -1.143494844436646
$data[11,7] = 4.552274703979492
$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = 9.923073837659718
$data[12,3] = -36.7092630550327
$data[12,4] = 23.39232449509959
$data[12,5] = -6.584865570068359
$data[12,6] = 3.639542579650879
$data[12,7] = 3.124820232391357
$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = -7.913778024561013
$data[13,3] = -1.389618308296037
$data[13,4] = 10.58147436029777
$data[13,5] = 0.5002082586288452
$data[13,6] = 12.75873851776123
$data[13,7] = -0.6872287392616272
$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = 22.76746183283226
$data[14,3] = -2.0813869372752
$data[14,4] = 18.41511002708894
$data[14,5] = -0.09793774783611291
$data[14,6] = -2.606635570526123
$data[14,7] = 1.64676570892334
$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = -10.96762119483142
$data[15,3] = -13.92066524471113
$data[15,4] = -14.51778084767945
$data[15,5] = 7.020159244537354
$data[15,6] = -8.503939628601074
$data[15,7] = -4.922720432281494
$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = 1.213289821849271
$data[16,3] = 1.019793061649235
$data[16,4] = -14.16178052565629
$data[16,5] = 4.105595588684082
$data[16,6] = -6.451707363128662
$data[16,7] = -2.353302240371704
$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = -7.05744662867424
$data[17,3] = -7.422974802250357
$data[17,4] = -4.164491389132067
$data[17,5] = -1.717672348022461
$data[17,6] = 1.564407467842102
$data[17,7] = 1.324523210525513
$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = 4.854760260603785
$data[18,3] = -2.305751679709482
$data[18,4] = 28.47659013605825
$data[18,5] = -12.60414218902588
$data[18,6] = -15.56717586517334
$data[18,7] = 2.47767186164856
$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = 3.906654929683194
$data[19,3] = 12.70785129771528
$data[19,4] = 6.405859976332387
$data[19,5] = 7.210841655731201
$data[19,6] = -12.07630348205566
$data[19,7] = 8.885769844055176
$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = 36.5629992182982
$data[20,3] = -53.32195605627643
$data[20,4] = 10.55735367024117
$data[20,5] = 1.924999475479126
$data[20,6] = -10.15509986877441
$data[20,7] = 0.8653942346572876
$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = 24.49468200994281
$data[21,3] = -0.5662195973802113
$data[21,4] = -16.59405546058915
$data[21,5] = 5.971939086914063
$data[21,6] = 14.02054500579834
$data[21,7] = -4.391685962677002
$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = -4.553085171799378
$data[22,3] = 1.362173399773924
$data[22,4] = -25.30839441912213
$data[22,5] = -2.465487957000732
$data[22,6] = 0.488490343093872
$data[22,7] = 5.172791004180908
$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = -26.14233835142689
$data[23,3] = -26.40326150402269
$data[23,4] = 16.80009495834494
$data[23,5] = -6.064484119415283
$data[23,6] = 1.030177354812622
$data[23,7] = 4.91339921951294
$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = -10.68485793376924
$data[24,3] = 1.987670397866637
$data[24,4] = 5.468654218302412
$data[24,5] = -0.608199417591095
$data[24,6] = 12.21864986419678
$data[24,7] = -3.169827461242676
$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = 11.52886270721615
$data[25,3] = 5.051807858825063
$data[25,4] = 21.48786570044184
$data[25,5] = -0.2609232068061828
$data[25,6] = 3.39666223526001
$data[25,7] = -2.653707027435303
$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = -43.52052723983955
$data[26,3] = -14.64106700646965
$data[26,4] = -27.88186194993908
$data[26,5] = 6.485929012298584
$data[26,6] = -1.756554484367371
$data[26,7] = -2.664892196655273
$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = -13.80793043283201
$data[27,3] = 3.756517190200164
$data[27,4] = -21.9789466857904
$data[27,5] = 1.498893618583679
$data[27,6] = -2.116081237792969
$data[27,7] = -2.421479225158691
$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -4.44925512771311
$data[28,3] = -2.411408648771987
$data[28,4] = -8.979007841774775
$data[28,5] = -0.8995492458343506
$data[28,6] = 1.140432238578796
$data[28,7] = 0.6262423396110535
$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = 4.063052345725026
$data[29,3] = -3.658945868996846
$data[29,4] = 20.40397767459713
$data[29,5] = -11.36417484283447
$data[29,6] = -11.00677871704102
$data[29,7] = -5.598630428314209

$ws.Range("A2:H31").Value = $data
